$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: add "Test Vectors" header (bold) in E14
$ws.Range("E14").Value = "Test Vectors"
$ws.Range("E14").Font.Bold = $true

# Row 15: E15
$ws.Range("E15").Value = "1. From idle S1->S2->S2"

# Row 16 (new row): E16
$ws.Range("E16").Value = "     first slide, then no more logo updates"

# Row 17: F17
$ws.Range("F17").Value = "*currently H->E->E works fine"

# Row 18: F18
$ws.Range("F18").Value = "*H->S->S works fine"

# Row 19: F19
$ws.Range("F19").Value = "*H->E->S is not"

# Row 20: G20
$ws.Range("G20").Value = "logo moves farther right & gets bigger? Do console prints to check"

# Update the selection in the sheet view
$ws.Range("G21").Select() | Out-Null
